$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell D1 from "是否为主持人" to "是否为昵称列"
$ws.Range("D1").Value = "是否为昵称列"

# Move the active cell/selection to D5 (matches saved cursor position in diff)
$ws.Range("D5").Select()
